$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before D ("Last Name" | new col | "Instrument"...),
#    shifting Instrument/Quantity and the trailing blank columns one to the right.
$ws.Columns("D").Insert()

# 2. Header row
$ws.Range("D1").Value = "Founder or Employee"

# 3. Row 2 (emp1) - fix the "Founder or Employee" import columns
$ws.Range("C2").Value = "Honest"
$ws.Range("D2").Value = "Founder"

# 4. Row 3 (emp2) - fix the "Founder or Employee" import columns
$ws.Range("C3").Value = "Good"
$ws.Range("D3").Value = "Founder"

# 5. New row 4 (emp3)
$ws.Range("A4").Value = "emp3@mycompany.com"
$ws.Range("B4").Value = "Emp3"
$ws.Range("C4").Value = "Awesome"
$ws.Range("D4").Value = "Employee"
$ws.Range("E4").Value = "Equity"
$ws.Range("F4").Value = 300

# 6. New row 5 (emp4)
$ws.Range("A5").Value = "emp4@mycompany.com"
$ws.Range("B5").Value = "Emp4"
$ws.Range("C5").Value = "Super"
$ws.Range("D5").Value = "Employee"
$ws.Range("E5").Value = "Preferred"
$ws.Range("F5").Value = 400

# 7. Hyperlinks for the two new email cells (matches the style already used by A2/A3)
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:emp3@mycompany.com", "", "", "emp3@mycompany.com")
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:emp4@mycompany.com", "", "", "emp4@mycompany.com")

# Re-apply the same look as the other hyperlinked email cells (A2/A3) since
# Hyperlinks.Add stamps its own built-in "Hyperlink" style.
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 8. Selection as left by the editor
[void]$ws.Range("F6").Select()
